# doctor_info.xlsx update:
#  1) Phone/Mobile columns (D,E) for the existing rows 2-4 were stored as
#     text; convert them to real numbers.
#  2) A new doctor row (row 5) is appended with the same column layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Convert existing Phone/Mobile cells to numeric values ---
$ws.Range("D2").Value = 2105149109
$ws.Range("E2").Value = 6970427541

$ws.Range("D3").Value = 2102110922
$ws.Range("E3").Value = 6977524030

$ws.Range("D4").Value = 2103455493
$ws.Range("E4").Value = 6932351230

# --- 2) Append new row 5 with the new doctor's information ---
$ws.Range("A5").Value = "ΜΠΑΚΑ Φ. ΦΩΤΕΙΝΗ"
$ws.Range("B5").Value = "Αργυροπούλου 11-13, Αθήνα - Κάτω Πατήσια, 11145, ΑΤΤΙΚΗΣ"
$ws.Range("C5").Value = "Δερματολόγος – Αφροδισιολόγος – Εφαρμογές Laser – Αισθητική Δερματολογία"

# Phone/Mobile for the new row stay as text (as typed from an import), so
# format the cells as Text before entering the digit strings - this keeps
# the leading zeros/format intact and avoids Excel auto-converting them to
# numbers.
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "2114001851"
$ws.Range("E5").Value = "6932431775"

# F5/G5 (Email / Ωρα) are empty for the new doctor, matching the other
# rows' empty Email/Hour cells. Typing a lone quote stores an actual empty
# text value in the cell instead of leaving it completely blank.
$ws.Range("F5").Value = "'"
$ws.Range("G5").Value = "'"
